# Add a new "time_taken" column (F) to the metadata sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: same text style as the other header cells (copy format from E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells F2/F3: plain text timestamps, no special formatting (like E2/E3).
$ws.Range("F2").Value = "2021-10-05 13:38:53.811580"
$ws.Range("F3").Value = "2021-10-05 13:38:53.811588"
